$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 4 updates
$ws.Range("N4").Value = 5
$ws.Range("O4").Value = 1.73
$ws.Range("P4").Value = 2
$ws.Range("Q4").Value = 2.44
$ws.Range("R4").Value = 1.56

# Row 5 updates
$ws.Range("G5").Value = 1.83
$ws.Range("H5").Value = 3.25
$ws.Range("I5").Value = 4.75
$ws.Range("J5").Value = 2.63
$ws.Range("K5").Value = 1.8
$ws.Range("L5").Value = 6.5
$ws.Range("W5").Value = 5.8
$ws.Range("AA5").Value = 2.67
$ws.Range("AB5").Value = 1.42
$ws.Range("AF5").Value = 15
$ws.Range("AK5").Value = 29
$ws.Range("AN5").Value = 7.5
$ws.Range("AO5").Value = 21
$ws.Range("AP5").Value = 19
$ws.Range("AS5").Value = 67

# Row 18 updates
$ws.Range("G18").Value = 1.25
$ws.Range("I18").Value = 10.25
$ws.Range("K18").Value = 2.57
$ws.Range("L18").Value = 8.25
$ws.Range("AA18").Value = 2.12
$ws.Range("AB18").Value = 1.57
$ws.Range("AC18").Value = 6.9
$ws.Range("AD18").Value = 5.8
$ws.Range("AF18").Value = 7.1
$ws.Range("AH18").Value = 32
$ws.Range("AJ18").Value = 10.75
$ws.Range("AK18").Value = 27
$ws.Range("AP18").Value = 32
